# Shooter wheel calibration: add a new "Friday 3/22 Evening" data/chart sheet (Sheet2)
# mirroring the existing Sheet1 calibration charts, and make Sheet2 the active tab.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 data -----------------------------------------------------------

$ws2.Range("A1").Value = "Friday 3/22 Evening"

$ws2.Range("A3").Value = "POWER"
$ws2.Range("B3").Value = "AVG_RPM"

$powerVals = @(0.25, 0.3, 0.35, 0.4, 0.45, 0.5, 0.55000000000000004, 0.6, 0.65, 0.7, 0.75, 0.8, 0.85, 0.9, 0.95, 1)
$rpmVals   = @(759.5, 1095.0999999999999, 1352.9, 1579.8, 1777.6, 1982.1, 2210.1, 2391.6999999999998, 2594.5, 2765, 2890.5, 3067.9, 3233.6, 3426.8, 3583, 3812.6)

for ($i = 0; $i -lt $powerVals.Length; $i++) {
    $row = 4 + $i
    $ws2.Cells.Item($row, 1).Value = $powerVals[$i]
    $ws2.Cells.Item($row, 2).Value = $rpmVals[$i]
}

$ws2.Range("A4:A19").NumberFormat = "0.00"
$ws2.Range("B4:B19").NumberFormat = "0.0"

# --- Chart 1 on Sheet2: AvgRpm-Power (full calibration curve) --------------

$co1 = $ws2.ChartObjects().Add(146.875, 36.75, 443.5, 216.0)
$chart1 = $co1.Chart
$chart1.ChartType = 74  # xlXYScatterLines

$ser1 = $chart1.SeriesCollection().NewSeries()
$ser1.Name = "AvgRpm-Power"
$ser1.XValues = "=Sheet2!`$A`$4:`$A`$19"
$ser1.Values = "=Sheet2!`$B`$4:`$B`$19"

$chart1.Axes(1).TickLabels.NumberFormat = "0.00"
$chart1.Axes(2).TickLabels.NumberFormat = "0.0"
$chart1.Axes(2).HasMajorGridlines = $true

$chart1.HasLegend = $true
$chart1.Legend.Position = -4152  # xlLegendPositionRight

# --- Chart 2 on Sheet2: mid range (linear trend of the mid-power band) -----

$co2 = $ws2.ChartObjects().Add(146.125, 281.25, 443.5, 216.0)
$chart2 = $co2.Chart
$chart2.ChartType = 74  # xlXYScatterLines

$ser2 = $chart2.SeriesCollection().NewSeries()
$ser2.Name = "mid range"
$ser2.XValues = "=Sheet2!`$B`$9:`$B`$15"
$ser2.Values = "=Sheet2!`$A`$9:`$A`$15"

$trend2 = $ser2.Trendlines().Add()
$trend2.Type = -4132  # xlLinear
$trend2.DisplayEquation = $true
$trend2.DisplayRSquared = $false

$chart2.Axes(1).TickLabels.NumberFormat = "0.0"
$chart2.Axes(2).TickLabels.NumberFormat = "0.00"
$chart2.Axes(2).HasMajorGridlines = $true

$chart2.HasLegend = $true
$chart2.Legend.Position = -4152  # xlLegendPositionRight

# --- Selection / active tab -------------------------------------------------

$ws2.Activate()
$ws2.Range("K23").Select()

Write-Output "edit complete"
